$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.80994075699194
$ws.Range("D2").Value = 4.630582415723802
$ws.Range("E2").Value = 19.76646936105121
$ws.Range("F2").Value = 22.00976897900639
$ws.Range("G2").Value = 24.82005699981466
$ws.Range("H2").Value = 13.34307326609354
$ws.Range("K2").Value = 11.429712908634
$ws.Range("L2").Value = 7.982111710387281
$ws.Range("O2").Value = 19.7759255541662
$ws.Range("B3").Value = 15.68849634596029
$ws.Range("D3").Value = 4.551244700977302
$ws.Range("E3").Value = 19.8643665441562
$ws.Range("F3").Value = 22.04504210212019
$ws.Range("G3").Value = 24.87744333548171
$ws.Range("H3").Value = 13.38928186310843
$ws.Range("K3").Value = 11.16674398827247
$ws.Range("L3").Value = 7.925738661152009
$ws.Range("O3").Value = 19.84876793321824
$ws.Range("B4").Value = 15.61651973169609
$ws.Range("D4").Value = 4.501198593786361
$ws.Range("E4").Value = 19.92810320630779
$ws.Range("F4").Value = 22.07297780406858
$ws.Range("G4").Value = 24.92170836421681
$ws.Range("H4").Value = 13.41984118329521
$ws.Range("K4").Value = 11.00087721554858
$ws.Range("L4").Value = 7.89209619456461
$ws.Range("O4").Value = 19.8980022803103
$ws.Range("B5").Value = 15.58786758209826
$ws.Range("D5").Value = 4.480483289318166
$ws.Range("E5").Value = 19.95498975408542
$ws.Range("F5").Value = 22.08593709605724
$ws.Range("G5").Value = 24.94200730057173
$ws.Range("H5").Value = 13.43284430854665
$ws.Range("K5").Value = 10.93224167799833
$ws.Range("L5").Value = 7.87864197738311
$ws.Range("O5").Value = 19.91919709670358
$ws.Range("B6").Value = 15.58315169775429
$ws.Range("D6").Value = 4.477024553222442
$ws.Range("E6").Value = 19.95950944423524
$ws.Range("F6").Value = 22.08818400693864
$ws.Range("G6").Value = 24.94551414163527
$ws.Range("H6").Value = 13.43503668359956
$ws.Range("K6").Value = 10.92078370937981
$ws.Range("L6").Value = 7.876423658640176
$ws.Range("O6").Value = 19.92278474436473
$ws.Range("B7").Value = 15.61613053453603
$ws.Range("D7").Value = 4.500920501274489
$ws.Range("E7").Value = 19.92846210799124
$ws.Range("F7").Value = 22.07314620438005
$ws.Range("G7").Value = 24.9219729838686
$ws.Range("H7").Value = 13.42001432124454
$ws.Range("K7").Value = 10.99995571085072
$ws.Range("L7").Value = 7.891913697783949
$ws.Range("O7").Value = 19.89828354312004
$ws.Range("B8").Value = 15.76754582095768
$ws.Range("D8").Value = 4.603511496606671
$ws.Range("E8").Value = 19.79947216525296
$ws.Range("F8").Value = 22.0206263763993
$ws.Range("G8").Value = 24.83796407589492
$ws.Range("H8").Value = 13.35855206921989
$ws.Range("K8").Value = 11.33999201645912
$ws.Range("L8").Value = 7.962479652760996
$ws.Range("O8").Value = 19.80010440101323
$ws.Range("B9").Value = 16.08377046230861
$ws.Range("D9").Value = 4.793539177866919
$ws.Range("E9").Value = 19.57525035106762
$ws.Range("F9").Value = 21.96756642930351
$ws.Range("G9").Value = 24.74525848361795
$ws.Range("H9").Value = 13.25537786268642
$ws.Range("K9").Value = 11.96927047426158
$ws.Range("L9").Value = 8.108051402719077
$ws.Range("O9").Value = 19.64345376284912
$ws.Range("B10").Value = 16.32612102763391
$ws.Range("D10").Value = 4.925622132587502
$ws.Range("E10").Value = 19.42795117425856
$ws.Range("F10").Value = 21.95913993290533
$ws.Range("G10").Value = 24.72150975579988
$ws.Range("H10").Value = 13.19015414134239
$ws.Range("K10").Value = 12.40552492856498
$ws.Range("L10").Value = 8.218682485309889
$ws.Range("O10").Value = 19.55037122621708
$ws.Range("B11").Value = 16.43816292684276
$ws.Range("D11").Value = 4.983933303797327
$ws.Range("E11").Value = 19.36471154904319
$ws.Range("F11").Value = 21.96195009170556
$ws.Range("G11").Value = 24.72040496207869
$ws.Range("H11").Value = 13.16277917305461
$ws.Range("K11").Value = 12.59772965756263
$ws.Range("L11").Value = 8.269655402530748
$ws.Range("O11").Value = 19.51283202762534
$ws.Range("B12").Value = 16.48081543911176
$ws.Range("D12").Value = 5.00574813367417
$ws.Range("E12").Value = 19.34130505640658
$ws.Range("F12").Value = 21.96396878928656
$ws.Range("G12").Value = 24.72138426596226
$ws.Range("H12").Value = 13.15274314775846
$ws.Range("K12").Value = 12.66956926742137
$ws.Range("L12").Value = 8.289037168405986
$ws.Range("O12").Value = 19.49930992579087
$ws.Range("B13").Value = 16.47162000012152
$ws.Range("D13").Value = 5.001061937168787
$ws.Range("E13").Value = 19.34632201514501
$ws.Range("F13").Value = 21.96349159417802
$ws.Range("G13").Value = 24.72111116252712
$ws.Range("H13").Value = 13.15488989327816
$ws.Range("K13").Value = 12.654139996673
$ws.Range("L13").Value = 8.284859641509348
$ws.Range("O13").Value = 19.50219128774994
$ws.Range("B14").Value = 16.44166766435767
$ws.Range("D14").Value = 4.985733422871593
$ws.Range("E14").Value = 19.36277504507867
$ws.Range("F14").Value = 21.96209704801944
$ws.Range("G14").Value = 24.72045750701757
$ws.Range("H14").Value = 13.16194688289758
$ws.Range("K14").Value = 12.60365911449324
$ws.Range("L14").Value = 8.271248437817933
$ws.Range("O14").Value = 19.51170565110731
$ws.Range("B15").Value = 16.42334925598167
$ws.Range("D15").Value = 4.976309257405374
$ws.Range("E15").Value = 19.37292343015259
$ws.Range("F15").Value = 21.96136712142676
$ws.Range("G15").Value = 24.72023920412089
$ws.Range("H15").Value = 13.1663125102857
$ws.Range("K15").Value = 12.57261387642048
$ws.Range("L15").Value = 8.262921124394053
$ws.Range("O15").Value = 19.51762381180247
$ws.Range("B16").Value = 16.31883209525099
$ws.Range("D16").Value = 4.921774685660531
$ws.Range("E16").Value = 19.43215977661002
$ws.Range("F16").Value = 21.9590899195668
$ws.Range("G16").Value = 24.72177741122684
$ws.Range("H16").Value = 13.1919893685405
$ws.Range("K16").Value = 12.39283416533471
$ws.Range("L16").Value = 8.21536313742765
$ws.Range("O16").Value = 19.55292137893705
$ws.Range("B17").Value = 16.25515028662812
$ws.Range("D17").Value = 4.887856903196754
$ws.Range("E17").Value = 19.46946373650324
$ws.Range("F17").Value = 21.95939420126383
$ws.Range("G17").Value = 24.72520774215831
$ws.Range("H17").Value = 13.20832937744391
$ws.Range("K17").Value = 12.2809108899882
$ws.Range("L17").Value = 8.186343799287361
$ws.Range("O17").Value = 19.57580751223963
$ws.Range("B18").Value = 16.21869321091276
$ws.Range("D18").Value = 4.868181940289378
$ws.Range("E18").Value = 19.49127464760493
$ws.Range("F18").Value = 21.96019462637307
$ws.Range("G18").Value = 24.72809357577156
$ws.Range("H18").Value = 13.21794378349986
$ws.Range("K18").Value = 12.21594980928314
$ws.Range("L18").Value = 8.169714289959556
$ws.Range("O18").Value = 19.58942310400479
$ws.Range("B19").Value = 16.2063798609039
$ws.Range("D19").Value = 4.861492118748758
$ws.Range("E19").Value = 19.49872038504557
$ws.Range("F19").Value = 21.96057306037693
$ws.Range("G19").Value = 24.72922732592653
$ws.Range("H19").Value = 13.22123616036227
$ws.Range("K19").Value = 12.19385594477114
$ws.Range("L19").Value = 8.16409482182714
$ws.Range("O19").Value = 19.59411068515145
$ws.Range("B20").Value = 16.26191186636067
$ws.Range("D20").Value = 4.891484809419445
$ws.Range("E20").Value = 19.46545596610411
$ws.Range("F20").Value = 21.95929708508963
$ws.Range("G20").Value = 24.72474807958326
$ws.Range("H20").Value = 13.2065675924239
$ws.Range("K20").Value = 12.29288628050237
$ws.Range("L20").Value = 8.189426674789392
$ws.Range("O20").Value = 19.57332443734566
$ws.Range("B21").Value = 16.45045956162589
$ws.Range("D21").Value = 4.990243093787706
$ws.Range("E21").Value = 19.35792771532518
$ws.Range("F21").Value = 21.96248076419869
$ws.Range("G21").Value = 24.72061155300927
$ws.Range("H21").Value = 13.1598651088697
$ws.Range("K21").Value = 12.6185125405865
$ws.Range("L21").Value = 8.275244330598664
$ws.Range("O21").Value = 19.5088922215779
$ws.Range("B22").Value = 16.57498033928835
$ws.Range("D22").Value = 5.053230150090671
$ws.Range("E22").Value = 19.2908047044743
$ws.Range("F22").Value = 21.97012436534384
$ws.Range("G22").Value = 24.72605503976493
$ws.Range("H22").Value = 13.13126754762651
$ws.Range("K22").Value = 12.82580647775083
$ws.Range("L22").Value = 8.331787797668062
$ws.Range("O22").Value = 19.47082339921008
$ws.Range("B23").Value = 16.50841396455301
$ws.Range("D23").Value = 5.019758791574315
$ws.Range("E23").Value = 19.32634127436701
$ws.Range("F23").Value = 21.96553630602734
$ws.Range("G23").Value = 24.72240371358582
$ws.Range("H23").Value = 13.1463543900875
$ws.Range("K23").Value = 12.7156890718682
$ws.Range("L23").Value = 8.301572177036585
$ws.Range("O23").Value = 19.49077093319639
$ws.Range("B24").Value = 16.25885447253584
$ws.Range("D24").Value = 4.889845177796688
$ws.Range("E24").Value = 19.46726674406952
$ws.Range("F24").Value = 21.95933904291709
$ws.Range("G24").Value = 24.72495304720665
$ws.Range("H24").Value = 13.20736340916537
$ws.Range("K24").Value = 12.287474118256
$ws.Range("L24").Value = 8.188032737255659
$ws.Range("O24").Value = 19.57444560888397
$ws.Range("B25").Value = 15.9963343342721
$ws.Range("D25").Value = 4.743403692561427
$ws.Range("E25").Value = 19.63284088335297
$ws.Range("F25").Value = 21.97655980144637
$ws.Range("G25").Value = 24.762574245607
$ws.Range("H25").Value = 13.28143152097683
$ws.Range("K25").Value = 11.80339751897476
$ws.Range("L25").Value = 8.067972136657529
$ws.Range("O25").Value = 19.68197646528248
